# Quarterly rollover: drop the oldest quarter (1399/06) and append the new
# quarter (1401/12). Every quarterly data column E:N shifts one column to
# the left and a freshly-reported value lands in the new N column. Also
# re-apply the (fixed) expense-classification ("read_price") split between
# "advertising" (row 13) and "other expenses" (row 19) for the 1401/06
# column, which nets to the same row total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N")

# Row 8 / Row 24: quarter-period header labels shown above the table.
$headers = @(
  "فصل سوم منتهی به 1399/09",
  "فصل چهارم منتهی به 1399/12",
  "فصل اول منتهی به 1400/03",
  "فصل دوم منتهی به 1400/06",
  "فصل سوم منتهی به 1400/09",
  "فصل چهارم منتهی به 1400/12",
  "فصل اول منتهی به 1401/03",
  "فصل دوم منتهی به 1401/06",
  "فصل سوم منتهی به 1401/09",
  "فصل چهارم منتهی به 1401/12"
)

# Data rows 10,11,12,13,14,15,16,17,18,19,20,26,27: new quarterly series
# (E:N) after the rollover + reclassification.
$rowData = @{
  10 = @(0,0,36059,53400,25927,37156,59578,36206,28901,150215)
  11 = @(0,0,0,0,0,0,0,0,0,0)
  12 = @(0,0,0,0,0,0,0,0,0,0)
  13 = @(2057,-1594,0,0,586,-586,0,0,0,0)
  14 = @(0,0,0,0,0,0,0,0,0,0)
  15 = @(181,264,340,563,703,519,573,594,442,1183)
  16 = @(139,128,230,263,265,325,352,544,549,680)
  17 = @(19343,21484,28929,33897,31631,33606,54802,38785,51946,57328)
  18 = @(0,0,0,0,0,0,0,0,0,0)
  19 = @(1689,3314,7182,-5497,3154,-3609,5978,3298,3119,546)
  20 = @(23409,23596,72740,82626,62266,67411,121283,79427,84957,209952)
  26 = @(55,58,58,58,57,57,55,55,53,55)
  27 = @(470,468,456,470,463,464,454,457,459,457)
}

for ($i = 0; $i -lt $cols.Length; $i++) {
  $col = $cols[$i]
  $ws.Range($col + "8").Value = $headers[$i]
  $ws.Range($col + "24").Value = $headers[$i]

  foreach ($r in $rowData.Keys) {
    $ws.Range($col + $r).Value = $rowData[$r][$i]
  }
}
